$wb = $excel.ActiveWorkbook

# --- BEPEfCT sheet: enable (set to 1) the carbon-tax exemption flags for
#     "agriculture and forestry" and "water and waste" sectors ---
$ws2 = $wb.Worksheets.Item("BEPEfCT")
$ws2.Range("B2").Value = 1
$ws2.Range("B25").Value = 1

# --- About sheet: append explanatory text describing the change ---
$ws1 = $wb.Worksheets.Item("About")
$ws1.Range("A13").Value = "In the U.S., we exempt agriculture and water and waste process emissions. Generally, "
$ws1.Range("A14").Value = "proposed taxes do not cover these sectors."

# --- Restore the on-screen selections seen in the saved workbook ---
$ws2.Range("B5").Select()
$ws1.Range("A15").Select()
